$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain text even if the new value looks numeric,
    # matching the inline-string (t="inlineStr") cells already in the sheet.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.683.98"
$ws.Range("E2").Value = "  +5.58%  "
Set-TextValue $ws.Range("D3") "3.520.88"
$ws.Range("E3").Value = "  +8.52%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "562.27"
$ws.Range("E5").Value = "  +8.59%  "
Set-TextValue $ws.Range("D6") "186.07"
$ws.Range("E6").Value = "  +10.82%  "
Set-TextValue $ws.Range("D7") "0.640"
$ws.Range("E7").Value = "  +10.52%  "
Set-TextValue $ws.Range("D8") "3.515.76"
$ws.Range("E8").Value = "  +8.34%  "
$ws.Range("E9").Value = "  +0.13%  "
Set-TextValue $ws.Range("D10") "0.644"
$ws.Range("E10").Value = "  +9.25%  "
Set-TextValue $ws.Range("D11") "0.156"
$ws.Range("E11").Value = "  +20.48%  "
Set-TextValue $ws.Range("D12") "55.70"
$ws.Range("E12").Value = "  +7.92%  "
Set-TextValue $ws.Range("D13") "0.0000280"
$ws.Range("E13").Value = "  +11.44%  "
Set-TextValue $ws.Range("D14") "9.45"
$ws.Range("E14").Value = "  +8.04%  "
Set-TextValue $ws.Range("D15") "4.070.81"
$ws.Range("E15").Value = "  +8.48%  "
Set-TextValue $ws.Range("D16") "3.512.55"
$ws.Range("E16").Value = "  +8.46%  "
Set-TextValue $ws.Range("D17") "18.71"
$ws.Range("E17").Value = "  +9.89%  "
$ws.Range("E18").Value = "  +6.64%  "
Set-TextValue $ws.Range("D19") "66.631.70"
$ws.Range("E19").Value = "  +5.93%  "
Set-TextValue $ws.Range("D20") "12.15"
$ws.Range("E20").Value = "  +11.18%  "
Set-TextValue $ws.Range("D21") "1.01"
$ws.Range("E21").Value = "  +7.67%  "
Set-TextValue $ws.Range("D22") "420.42"
$ws.Range("E22").Value = "  +13.51%  "
Set-TextValue $ws.Range("D23") "4.11"
$ws.Range("E23").Value = "  +14.04%  "
Set-TextValue $ws.Range("D24") "86.24"
$ws.Range("E24").Value = "  +7.62%  "
$ws.Range("E25").Value = "  +1.81%  "
Set-TextValue $ws.Range("D26") "11.01"
$ws.Range("E26").Value = "  -0.06%  "
Set-TextValue $ws.Range("D27") "2.93"
$ws.Range("E27").Value = "  +11.37%  "
Set-TextValue $ws.Range("D28") "12.45"
$ws.Range("E28").Value = "  +13.45%  "
Set-TextValue $ws.Range("D29") "6.12"
$ws.Range("E29").Value = "  -0.03%  "
Set-TextValue $ws.Range("D30") "9.25"
$ws.Range("E30").Value = "  +16.88%  "
Set-TextValue $ws.Range("D31") "30.52"
$ws.Range("E31").Value = "  +9.32%  "
Set-TextValue $ws.Range("D32") "6.75"
$ws.Range("E32").Value = "  +4.19%  "
Set-TextValue $ws.Range("D33") "619.17"
$ws.Range("E33").Value = "  +0.58%  "
Set-TextValue $ws.Range("D34") "11.95"
$ws.Range("E34").Value = "  +8.99%  "
$ws.Range("E35").Value = "  +9.33%  "
Set-TextValue $ws.Range("D36") "60.50"
$ws.Range("E36").Value = "  +7.38%  "
Set-TextValue $ws.Range("D37") "0.151"
$ws.Range("E37").Value = "  +23.82%  "
Set-TextValue $ws.Range("D38") "0.0₃0821"
$ws.Range("E38").Value = "  +16.75%  "
Set-TextValue $ws.Range("D39") "38.40"
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("E40").Value = "  -0.07%  "
Set-TextValue $ws.Range("D41") "0.389"
$ws.Range("E41").Value = "  +5.47%  "
Set-TextValue $ws.Range("D42") "3.38"
$ws.Range("E42").Value = "  +10.98%  "
Set-TextValue $ws.Range("D43") "3.128.94"
$ws.Range("E43").Value = "  +11.35%  "
Set-TextValue $ws.Range("D44") "0.997"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  +13.78%  "
Set-TextValue $ws.Range("D47") "3.30"
$ws.Range("E47").Value = "  +9.93%  "
Set-TextValue $ws.Range("D48") "0.0420"
$ws.Range("E48").Value = "  +8.35%  "
Set-TextValue $ws.Range("D49") "2.73"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("E50").Value = "  +9.49%  "
Set-TextValue $ws.Range("D51") "139.94"
$ws.Range("E51").Value = "  +2.56%  "
